$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "item" formula label cell to use the localized label/display-values helper
$ws.Range("A6").Value = '${item.getLocalizedLabelAndDisplayValues(locale)}'

# B5 text remains "${cityId.value}" (string table was reshuffled but displayed text is unchanged)
$ws.Range("B5").Value = '${cityId.value}'

# Update the selected/active cell from F12 to A7
$ws.Range("A7").Select()

$wb.Save()
